$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 12293
$ws1.Range("F14").Value = 13763
$ws1.Range("F22").Value = 4842
$ws1.Range("F23").Value = 210

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 12293
$ws4.Range("F14").Value = 13763
$ws4.Range("F22").Value = 4842
$ws4.Range("F23").Value = 210
